# Updated cryptos list (price/volume refresh + Mantle/VeChain row swap).
# Note: several "Price" values look like plain numbers (e.g. 576.35).
# Column D is plain text in this sheet, so those are written with a
# leading apostrophe to force text entry (avoiding numeric auto-conversion),
# then the cell style is reset to "Normal" so no extra number-format / style
# ends up attached to the cell (matching the original formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.085.77'
$ws.Range("E2").Value = '  +0.52%  '

$ws.Range("D3").Value = '2.760.97'
$ws.Range("E3").Value = '  +1.12%  '

$ws.Range("E4").Value = '  -0.35%  '

$ws.Range("D5").Value = '''576.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.74%  '

$ws.Range("D6").Value = '''159.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.99%  '

$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("E8").Value = '  -3.40%  '

$ws.Range("E9").Value = '  -1.71%  '

$ws.Range("D10").Value = '''0.164'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.13%  '

$ws.Range("D11").Value = '''0.386'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.00%  '

$ws.Range("D12").Value = '''5.67'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -16.21%  '

$ws.Range("D13").Value = '3.249.45'
$ws.Range("E13").Value = '  +0.56%  '

$ws.Range("D14").Value = '''26.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.29%  '

$ws.Range("D15").Value = '63.841.66'
$ws.Range("E15").Value = '  +0.15%  '

$ws.Range("E16").Value = '  -2.25%  '

$ws.Range("D17").Value = '2.763.67'
$ws.Range("E17").Value = '  +0.21%  '

$ws.Range("D18").Value = '''12.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.80%  '

$ws.Range("D19").Value = '''4.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.55%  '

$ws.Range("D20").Value = '''357.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.05%  '

$ws.Range("D21").Value = '''6.75'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.75%  '

$ws.Range("E22").Value = '  +0.19%  '

$ws.Range("D23").Value = '''0.533'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.43%  '

$ws.Range("D24").Value = '''65.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.29%  '

$ws.Range("E25").Value = '  -0.63%  '

$ws.Range("D26").Value = '''8.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("D28").Value = '0.0₃0911'
$ws.Range("E28").Value = '  -1.08%  '

$ws.Range("E29").Value = '  -2.90%  '

$ws.Range("D30").Value = '''7.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.20%  '

$ws.Range("E31").Value = '  -0.52%  '

$ws.Range("D32").Value = '''169.70'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.27%  '

$ws.Range("D33").Value = '''20.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.76%  '

$ws.Range("E34").Value = '  -0.19%  '

$ws.Range("D35").Value = '''1.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.83%  '

$ws.Range("D36").Value = '''0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("D37").Value = '''1.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.76%  '

$ws.Range("D38").Value = '''1.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.27%  '

$ws.Range("D39").Value = '''6.31'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.41%  '

$ws.Range("D40").Value = '''339.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.26%  '

$ws.Range("E41").Value = '  -2.28%  '

$ws.Range("E42").Value = '  -0.63%  '

$ws.Range("D43").Value = '''21.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.50%  '

$ws.Range("D44").Value = '''21.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.38%  '

$ws.Range("D45").Value = '''0.0591'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.67%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '''0.635'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.97%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '''0.0256'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.71%  '

$ws.Range("E48").Value = '  -0.55%  '

$ws.Range("D49").Value = '''135.95'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.51%  '

$ws.Range("D50").Value = '''0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.31%  '

$ws.Range("D51").Value = '''11.07'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.13%  '
